# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.971.88'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '2.636.02'
$ws.Range('E3').Value = '  +2.47%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '''513.95'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '''143.92'
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('D7').Value = '''0.995'
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('D8').Value = '''0.567'
$ws.Range('E8').Value = '  +1.97%  '
$ws.Range('D9').Value = '2.664.26'
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('D10').Value = '''6.30'
$ws.Range('E10').Value = '  +1.64%  '
$ws.Range('E11').Value = '  +3.03%  '
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').Value = '''0.126'
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').Value = '3.096.94'
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('D15').Value = '58.946.63'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '''21.08'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('D18').Value = '2.654.36'
$ws.Range('E18').Value = '  +2.62%  '
$ws.Range('D19').Value = '''4.54'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '''341.23'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').Value = '''10.40'
$ws.Range('E21').Value = '  +2.65%  '
$ws.Range('D22').Value = '''6.10'
$ws.Range('E22').Value = '  +2.04%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = '''60.97'
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('D25').Value = '''0.420'
$ws.Range('E25').Value = '  +2.50%  '
$ws.Range('D26').Value = '2.752.99'
$ws.Range('E26').Value = '  +1.50%  '
$ws.Range('D27').Value = '''0.993'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').Value = '''0.161'
$ws.Range('E28').Value = '  +3.55%  '
$ws.Range('D29').Value = '0.0₃0803'
$ws.Range('E29').Value = '  +3.11%  '
$ws.Range('D30').Value = '''7.11'
$ws.Range('E30').Value = '  +3.39%  '
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').Value = '''6.38'
$ws.Range('E32').Value = '  +8.83%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''18.90'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').Value = '''1.58'
$ws.Range('E34').Value = '  +1.89%  '
$ws.Range('D35').Value = '''148.89'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  +14.14%  '
$ws.Range('E37').Value = '  +4.56%  '
$ws.Range('E38').Value = '  +3.36%  '
$ws.Range('D39').Value = '''0.854'
$ws.Range('E39').Value = '  +3.62%  '
$ws.Range('D40').Value = '''36.57'
$ws.Range('E40').Value = '  +1.28%  '
$ws.Range('E41').Value = '  +4.07%  '
$ws.Range('E42').Value = '  +1.33%  '
$ws.Range('D43').Value = '''281.89'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = '''0.615'
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('D46').Value = '''0.0981'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '''19.48'
$ws.Range('E47').Value = '  +4.69%  '
$ws.Range('D48').Value = '''0.0532'
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '''0.0230'
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').Value = '''4.71'
$ws.Range('E50').Value = '  +4.60%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '''10.28'
$ws.Range('E51').Value = '  -0.59%  '
